# Apply edit described by commit "Identified papers which i have access to"
#
# Summary of the change:
#  - On "Initial Search", four papers (rows 4, 12, 15, 38 in column D) are
#    marked red (same highlight already used for other excluded papers).
#  - Those same four papers are removed from the "Records Sought for
#    retrieval" list, since they no longer need to be sought.
#  - The remaining papers in "Records Sought for retrieval" get their B
#    column coloured red or green depending on whether access was found.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Initial Search")
$ws2 = $wb.Worksheets.Item("Records Sought for retrieval")

# Colours used elsewhere in the workbook (OLE_COLOR = BGR order):
#   red   -> RGB(FF0000) -> OLE 0x0000FF -> 255
#   green -> RGB(00B050) -> OLE 0x50B000 -> 5287936
$RED   = 255
$GREEN = 5287936

# --- 1. Mark the four inaccessible papers red on "Initial Search" ---
$ws1.Range("D4").Interior.Color  = $RED
$ws1.Range("D12").Interior.Color = $RED
$ws1.Range("D15").Interior.Color = $RED
$ws1.Range("D38").Interior.Color = $RED

# --- 2. Remove those four papers from "Records Sought for retrieval" ---
# (original rows 3, 9, 12, 24 -- delete bottom-up so row numbers stay valid)
$ws2.Rows("24").Delete()
$ws2.Rows("12").Delete()
$ws2.Rows("9").Delete()
$ws2.Rows("3").Delete()

# --- 3. Colour the B column for the remaining 27 papers (rows 2-28) ---
$styles = @($GREEN, $GREEN, $RED, $RED, $RED, $GREEN, $GREEN, $GREEN, $GREEN, $GREEN, $GREEN, $RED, $GREEN, $GREEN, $GREEN, $GREEN, $RED, $GREEN, $GREEN, $GREEN, $RED, $GREEN, $GREEN, $GREEN, $GREEN, $GREEN, $GREEN)

for ($i = 0; $i -lt $styles.Length; $i++) {
    $row = $i + 2
    $ws2.Range("B$row").Interior.Color = $styles[$i]
}

# --- 4. Best-effort view state to match the saved workbook ---
# "Initial Search" keeps its own selection (D38 ends up selected there,
# matching the cell that was just marked red), while "Records Sought for
# retrieval" remains the active tab with A11 selected.
$ws1.Activate() | Out-Null
$ws1.Range("D38").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A11").Select() | Out-Null
